$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.291.68'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.776.66'
$ws.Range('E3').Value = '  +3.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.34'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.39%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5199'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +8.97%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3700'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.84'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07398'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.090'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.39%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('E13').Value = '  +3.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.073'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = '1.767.07'
$ws.Range('E15').Value = '  +2.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.969'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.20'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001048'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06439'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('E21').Value = '  +1.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.822'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.78%  '
$ws.Range('D23').Value = '27.320.83'
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('E24').Value = '  +4.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.119'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.07%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.327'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +11.11%  '
$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D29').Value = '1.972.29'
$ws.Range('E29').Value = '  +3.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '121.34'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.064'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09794'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.579'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.625'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02245'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05979'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.26%  '
$ws.Range('E37').Value = '  +1.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.843'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6144'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.68%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2020'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.432'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.103'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.139'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.15'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5771'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.630'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.890'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.95%  '
$ws.Range('E49').Value = '  +2.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06714'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.49%  '
